# Fruta / hortaliza, semanal
# Inserts a new weekly price record as row 144 (shifting the existing
# rows 144-177 down to 145-178), matching the new "Acelga" observation
# for Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by one row, starting at row 144.
$ws.Rows.Item(144).Insert()

# Populate the newly opened row 144 with the new weekly observation.
$ws.Cells.Item(144, 1).Value  = 7
$ws.Cells.Item(144, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(144, 3).Value  = "Ñuble"
$ws.Cells.Item(144, 4).Value  = 44508
$ws.Cells.Item(144, 5).Value  = 16
$ws.Cells.Item(144, 6).Value  = 100112009
$ws.Cells.Item(144, 7).Value  = "Acelga"
$ws.Cells.Item(144, 8).Value  = "Sin especificar"
$ws.Cells.Item(144, 9).Value  = "Primera"
$ws.Cells.Item(144, 10).Value = 120
$ws.Cells.Item(144, 11).Value = 350
$ws.Cells.Item(144, 12).Value = 400
$ws.Cells.Item(144, 13).Value = 375
$ws.Cells.Item(144, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(144, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(144, 16).Value = 375
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = "Hortaliza"
